$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed values from the
# scheduled GitHub Actions symbol-list sync. Values are written with a leading
# apostrophe so Excel keeps storing them as text (matching the sheet's existing
# inline-string cells) instead of re-typing them as numbers/percentages.
$ws.Range("D2").Value = "'326.44"
$ws.Range("E2").Value = "'-1.26%"
$ws.Range("D3").Value = "'45.15"
$ws.Range("E3").Value = "'2.91%"
$ws.Range("D4").Value = "'5.566"
$ws.Range("E4").Value = "'-6.79%"
$ws.Range("D5").Value = "'0.08092"
$ws.Range("E5").Value = "'-2.47%"
$ws.Range("D6").Value = "'8.710"
$ws.Range("E6").Value = "'-0.92%"
$ws.Range("E7").Value = "'-3.66%"
$ws.Range("D8").Value = "'1.909"
$ws.Range("E8").Value = "'-2.75%"
$ws.Range("E9").Value = "'-6.20%"
$ws.Range("D10").Value = "'0.9495"
$ws.Range("E10").Value = "'2.05%"
$ws.Range("D11").Value = "'0.1163"
$ws.Range("E11").Value = "'-6.87%"
$ws.Range("D12").Value = "'0.1893"
$ws.Range("E12").Value = "'-3.32%"
$ws.Range("D13").Value = "'0.1017"
$ws.Range("E13").Value = "'6.30%"
$ws.Range("D14").Value = "'0.04137"
$ws.Range("E14").Value = "'4.05%"
$ws.Range("E15").Value = "'-0.13%"
$ws.Range("D16").Value = "'0.001274"
$ws.Range("E16").Value = "'-2.52%"
$ws.Range("D17").Value = "'0.006059"
$ws.Range("E17").Value = "'0.67%"
$ws.Range("D18").Value = "'3.620"
$ws.Range("E18").Value = "'2.45%"
$ws.Range("E19").Value = "'-0.68%"
$ws.Range("D20").Value = "'8.535"
$ws.Range("E20").Value = "'-6.30%"
$ws.Range("D21").Value = "'0.1383"
$ws.Range("E21").Value = "'0.66%"
$ws.Range("E22").Value = "'2.40%"
$ws.Range("E23").Value = "'-3.40%"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("E24").Value = "'-0.53%"
$ws.Range("D25").Value = "'0.004594"
$ws.Range("E25").Value = "'4.47%"
$ws.Range("E26").Value = "'3.43%"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("D39").Value = "'0.02690"
$ws.Range("E39").Value = "'-4.68%"
$ws.Range("D40").Value = "'0.05565"
$ws.Range("E40").Value = "'-1.04%"
$ws.Range("E41").Value = "'25.10%"
$ws.Range("D42").Value = "'0.007706"
$ws.Range("E42").Value = "'-2.78%"
$ws.Range("D43").Value = "'0.1395"
$ws.Range("E43").Value = "'-1.91%"
$ws.Range("D44").Value = "'0.002076"
$ws.Range("E44").Value = "'-1.86%"
$ws.Range("D45").Value = "'0.008690"
$ws.Range("E45").Value = "'-0.75%"
$ws.Range("D46").Value = "'0.00007131"
$ws.Range("E46").Value = "'-2.51%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D48").Value = "'0.003509"
$ws.Range("E48").Value = "'-2.64%"
$ws.Range("D49").Value = "'0.002277"
$ws.Range("E49").Value = "'-0.26%"
$ws.Range("D51").Value = "'0.0002006"
